$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5
$ws.Range("I4").Value = 6.5
$ws.Range("L4").Value = 7.5
$ws.Range("Z4").Value = 9.5
$ws.Range("AF4").Value = 101
$ws.Range("AG4").Value = 13
$ws.Range("AH4").Value = 34
$ws.Range("AL4").Value = 67
$ws.Range("AO4").Value = 7.5
$ws.Range("AZ4").Value = 201
$ws.Range("BA4").Value = 251
$ws.Range("M10").Value = 1.11
$ws.Range("N10").Value = 6.5
$ws.Range("M14").Value = 1.06
$ws.Range("N14").Value = 9.5
$ws.Range("O14").Value = 1.3
$ws.Range("P14").Value = 3.4
$ws.Range("Q14").Value = 2.03
$ws.Range("R14").Value = 1.83
$ws.Range("S14").Value = 1.4
$ws.Range("T14").Value = 2.75
$ws.Range("AC14").Value = 9.5
$ws.Range("AG14").Value = 6.5
$ws.Range("AK14").Value = 13
$ws.Range("AS14").Value = 301
$ws.Range("AT14").Value = 2.75
$ws.Range("AX14").Value = 8
$ws.Range("Q15").Value = 2.5
$ws.Range("R15").Value = 1.5
$ws.Range("G18").Value = 7.3
$ws.Range("H18").Value = 5
$ws.Range("J18").Value = 5.9
$ws.Range("K18").Value = 2.8
$ws.Range("P18").Value = 6.2
$ws.Range("Q18").Value = 1.3
$ws.Range("R18").Value = 3.25
$ws.Range("S18").Value = 1.18
$ws.Range("T18").Value = 4.3
$ws.Range("U18").Value = 1.47
$ws.Range("V18").Value = 2.5
$ws.Range("W18").Value = 37
$ws.Range("X18").Value = 70
$ws.Range("Y18").Value = 23
$ws.Range("Z18").Value = 175
$ws.Range("AA18").Value = 65
$ws.Range("AC18").Value = 10.75
$ws.Range("AD18").Value = 11.75
$ws.Range("AG18").Value = 13.5
$ws.Range("AH18").Value = 10.25
$ws.Range("AI18").Value = 9
$ws.Range("AJ18").Value = 11.25
$ws.Range("AL18").Value = 16
$ws.Range("AO18").Value = 35
$ws.Range("AQ18").Value = 175
$ws.Range("AT18").Value = 4.3
$ws.Range("AU18").Value = 6.9
$ws.Range("AW18").Value = 3.9
$ws.Range("AY18").Value = 10.75
$ws.Range("AZ18").Value = 13.5
$ws.Range("BA18").Value = 24
$ws.Range("BC18").Value = 350
